$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Slit2"
$ws.Range("C2").Value2 = "Sdc1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.143896
$ws.Range("H2").Value2 = 0.431688
$ws.Range("I2").Value2 = 0.02807111181859822
$ws.Range("J2").Value2 = 0.02807111181859822
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.9159313333333333
$ws.Range("N2").Value2 = 2.747794
$ws.Range("O2").Value2 = 0.03641350786393945
$ws.Range("P2").Value2 = 0.03641350786393944
$ws.Range("Q2").Value2 = 0.1317988551413333
$ws.Range("R2").Value2 = 1.186189696272
$ws.Range("S2").Value2 = 0.00102216765095605
$ws.Range("T2").Value2 = 0.00102216765095605

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Slit2"
$ws.Range("C3").Value2 = "Sdc1"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.143896
$ws.Range("H3").Value2 = 0.431688
$ws.Range("I3").Value2 = 0.02807111181859822
$ws.Range("J3").Value2 = 0.02807111181859822
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 10.56834633333333
$ws.Range("N3").Value2 = 31.705039
$ws.Range("O3").Value2 = 0.4201521973455823
$ws.Range("P3").Value2 = 0.4201521973455822
$ws.Range("Q3").Value2 = 1.520742763981333
$ws.Range("R3").Value2 = 13.686684875832
$ws.Range("S3").Value2 = 0.01179413931251759
$ws.Range("T3").Value2 = 0.01179413931251759

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Slit2"
$ws.Range("C4").Value2 = "Sdc1"
$ws.Range("D4").Value2 = "M2"
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.143896
$ws.Range("H4").Value2 = 0.431688
$ws.Range("I4").Value2 = 0.02807111181859822
$ws.Range("J4").Value2 = 0.02807111181859822
$ws.Range("K4").Value2 = 2
$ws.Range("L4").Value2 = 0.6666666666666666
$ws.Range("M4").Value2 = 1.220967333333333
$ws.Range("N4").Value2 = 3.662902
$ws.Range("O4").Value2 = 0.04854043308262539
$ws.Range("P4").Value2 = 0.04854043308262537
$ws.Range("Q4").Value2 = 0.1756923153973333
$ws.Range("R4").Value2 = 1.581230838576
$ws.Range("S4").Value2 = 0.001362583924785562
$ws.Range("T4").Value2 = 0.001362583924785561

$ws.Range("A5").Value2 = "ECs"
$ws.Range("B5").Value2 = "Slit2"
$ws.Range("C5").Value2 = "Sdc1"
$ws.Range("D5").Value2 = "sCs"
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.143896
$ws.Range("H5").Value2 = 0.431688
$ws.Range("I5").Value2 = 0.02807111181859822
$ws.Range("J5").Value2 = 0.02807111181859822
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 12.44836933333333
$ws.Range("N5").Value2 = 37.345108
$ws.Range("O5").Value2 = 0.494893861707853
$ws.Range("P5").Value2 = 0.4948938617078529
$ws.Range("Q5").Value2 = 1.791270553589333
$ws.Range("R5").Value2 = 16.121434982304
$ws.Range("S5").Value2 = 0.01389222093033903
$ws.Range("T5").Value2 = 0.01389222093033902

$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Slit2"
$ws.Range("C6").Value2 = "Sdc1"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 4.277274333333334
$ws.Range("H6").Value2 = 12.831823
$ws.Range("I6").Value2 = 0.8344071140950421
$ws.Range("J6").Value2 = 0.8344071140950421
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.9159313333333333
$ws.Range("N6").Value2 = 2.747794
$ws.Range("O6").Value2 = 0.03641350786393945
$ws.Range("P6").Value2 = 0.03641350786393944
$ws.Range("Q6").Value2 = 3.917689583162445
$ws.Range("R6").Value2 = 35.259206248462
$ws.Range("S6").Value2 = 0.03038369001082684
$ws.Range("T6").Value2 = 0.03038369001082683

$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Slit2"
$ws.Range("C7").Value2 = "Sdc1"
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 4.277274333333334
$ws.Range("H7").Value2 = 12.831823
$ws.Range("I7").Value2 = 0.8344071140950421
$ws.Range("J7").Value2 = 0.8344071140950421
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 10.56834633333333
$ws.Range("N7").Value2 = 31.705039
$ws.Range("O7").Value2 = 0.4201521973455823
$ws.Range("P7").Value2 = 0.4201521973455822
$ws.Range("Q7").Value2 = 45.20371651734411
$ws.Range("R7").Value2 = 406.833448656097
$ws.Range("S7").Value2 = 0.3505779824678179
$ws.Range("T7").Value2 = 0.3505779824678179

$ws.Range("A8").Value2 = "FAPs"
$ws.Range("B8").Value2 = "Slit2"
$ws.Range("C8").Value2 = "Sdc1"
$ws.Range("D8").Value2 = "M2"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 4.277274333333334
$ws.Range("H8").Value2 = 12.831823
$ws.Range("I8").Value2 = 0.8344071140950421
$ws.Range("J8").Value2 = 0.8344071140950421
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 1.220967333333333
$ws.Range("N8").Value2 = 3.662902
$ws.Range("O8").Value2 = 0.04854043308262539
$ws.Range("P8").Value2 = 0.04854043308262537
$ws.Range("Q8").Value2 = 5.222412236705111
$ws.Range("R8").Value2 = 47.001710130346
$ws.Range("S8").Value2 = 0.04050248268539695
$ws.Range("T8").Value2 = 0.04050248268539695

$ws.Range("A9").Value2 = "FAPs"
$ws.Range("B9").Value2 = "Slit2"
$ws.Range("C9").Value2 = "Sdc1"
$ws.Range("D9").Value2 = "sCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 4.277274333333334
$ws.Range("H9").Value2 = 12.831823
$ws.Range("I9").Value2 = 0.8344071140950421
$ws.Range("J9").Value2 = 0.8344071140950421
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 12.44836933333333
$ws.Range("N9").Value2 = 37.345108
$ws.Range("O9").Value2 = 0.494893861707853
$ws.Range("P9").Value2 = 0.4948938617078529
$ws.Range("Q9").Value2 = 53.24509064132045
$ws.Range("R9").Value2 = 479.205815771884
$ws.Range("S9").Value2 = 0.4129429589310005
$ws.Range("T9").Value2 = 0.4129429589310004

$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Slit2"
$ws.Range("C10").Value2 = "Sdc1"
$ws.Range("D10").Value2 = "ECs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 0.7049536666666666
$ws.Range("H10").Value2 = 2.114861
$ws.Range("I10").Value2 = 0.1375217740863597
$ws.Range("J10").Value2 = 0.1375217740863597
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.9159313333333333
$ws.Range("N10").Value2 = 2.747794
$ws.Range("O10").Value2 = 0.03641350786393945
$ws.Range("P10").Value2 = 0.03641350786393944
$ws.Range("Q10").Value2 = 0.6456891518482222
$ws.Range("R10").Value2 = 5.811202366633999
$ws.Range("S10").Value2 = 0.005007650202156564
$ws.Range("T10").Value2 = 0.005007650202156563

$ws.Range("A11").Value2 = "sCs"
$ws.Range("B11").Value2 = "Slit2"
$ws.Range("C11").Value2 = "Sdc1"
$ws.Range("D11").Value2 = "FAPs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 0.7049536666666666
$ws.Range("H11").Value2 = 2.114861
$ws.Range("I11").Value2 = 0.1375217740863597
$ws.Range("J11").Value2 = 0.1375217740863597
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 10.56834633333333
$ws.Range("N11").Value2 = 31.705039
$ws.Range("O11").Value2 = 0.4201521973455823
$ws.Range("P11").Value2 = 0.4201521973455822
$ws.Range("Q11").Value2 = 7.450194498286554
$ws.Range("R11").Value2 = 67.05175048457899
$ws.Range("S11").Value2 = 0.05778007556524679
$ws.Range("T11").Value2 = 0.05778007556524679

$ws.Range("A12").Value2 = "sCs"
$ws.Range("B12").Value2 = "Slit2"
$ws.Range("C12").Value2 = "Sdc1"
$ws.Range("D12").Value2 = "M2"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 0.7049536666666666
$ws.Range("H12").Value2 = 2.114861
$ws.Range("I12").Value2 = 0.1375217740863597
$ws.Range("J12").Value2 = 0.1375217740863597
$ws.Range("K12").Value2 = 2
$ws.Range("L12").Value2 = 0.6666666666666666
$ws.Range("M12").Value2 = 1.220967333333333
$ws.Range("N12").Value2 = 3.662902
$ws.Range("O12").Value2 = 0.04854043308262539
$ws.Range("P12").Value2 = 0.04854043308262537
$ws.Range("Q12").Value2 = 0.8607253985135555
$ws.Range("R12").Value2 = 7.746528586622
$ws.Range("S12").Value2 = 0.00667536647244287
$ws.Range("T12").Value2 = 0.006675366472442869

$ws.Range("A13").Value2 = "sCs"
$ws.Range("B13").Value2 = "Slit2"
$ws.Range("C13").Value2 = "Sdc1"
$ws.Range("D13").Value2 = "sCs"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 0.7049536666666666
$ws.Range("H13").Value2 = 2.114861
$ws.Range("I13").Value2 = 0.1375217740863597
$ws.Range("J13").Value2 = 0.1375217740863597
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 12.44836933333333
$ws.Range("N13").Value2 = 37.345108
$ws.Range("O13").Value2 = 0.494893861707853
$ws.Range("P13").Value2 = 0.4948938617078529
$ws.Range("Q13").Value2 = 8.775523605554222
$ws.Range("R13").Value2 = 78.979712449988
$ws.Range("S13").Value2 = 0.06805868184651351
$ws.Range("T13").Value2 = 0.06805868184651351
